# fayoumi (add course module)
# Update the Name and Phone_Number values on the data row (row 2):
#   - Name (A2): "Test Automation Fayoumi31" -> "Mohd122 Test3 Automation Fayoumi31"
#   - Phone_Number (B2): numeric 791000071 -> text "792101070" (entered with a leading
#     apostrophe, as a real user would, so Excel stores it as text with quotePrefix)
# Then move the active selection to B4, as reflected in the sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mohd122 Test3 Automation Fayoumi31"
$ws.Range("B2").Value = "'792101070"

$ws.Range("B4").Select()
